$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04036098626811
$ws.Range("D2").Value = 1.045323154005079
$ws.Range("E2").Value = 1.044017075773781
$ws.Range("F2").Value = 1.054369798885633
$ws.Range("I2").Value = 1.037319484472095
$ws.Range("J2").Value = 1.045448009521716
$ws.Range("K2").Value = 1.048091716843552
$ws.Range("L2").Value = 1.046789310677181
$ws.Range("M2").Value = 1.057113197572296
$ws.Range("N2").Value = 1.018938139772754

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041664731200635
$ws.Range("D3").Value = 1.046295525882199
$ws.Range("E3").Value = 1.045270203998451
$ws.Range("F3").Value = 1.055487872252348
$ws.Range("I3").Value = 1.037581422789206
$ws.Range("J3").Value = 1.046395572104535
$ws.Range("K3").Value = 1.048875086737022
$ws.Range("L3").Value = 1.04785243758687
$ws.Range("M3").Value = 1.058043721382023
$ws.Range("N3").Value = 1.01926265019584

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042507907013803
$ws.Range("D4").Value = 1.046924163470993
$ws.Range("E4").Value = 1.046080991296465
$ws.Range("F4").Value = 1.056211065048608
$ws.Range("I4").Value = 1.037749336771969
$ws.Range("J4").Value = 1.047007808971803
$ws.Range("K4").Value = 1.049380812403989
$ws.Range("L4").Value = 1.048539726702988
$ws.Range("M4").Value = 1.05864497005411
$ws.Range("N4").Value = 1.019472075284093

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042862277363653
$ws.Range("D5").Value = 1.047188311979421
$ws.Range("E5").Value = 1.046421831873959
$ws.Range("F5").Value = 1.056515031065628
$ws.Range("I5").Value = 1.037819550638078
$ws.Range("J5").Value = 1.047264979983705
$ws.Range("K5").Value = 1.04959314151071
$ws.Range("L5").Value = 1.048828515571004
$ws.Range("M5").Value = 1.058897530100245
$ws.Range("N5").Value = 1.01955998530304

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042921771889431
$ws.Range("D6").Value = 1.047232656051278
$ws.Range("E6").Value = 1.04647905973366
$ws.Range("F6").Value = 1.0565660645801
$ws.Range("I6").Value = 1.037831317752486
$ws.Range("J6").Value = 1.047308147673457
$ws.Range("K6").Value = 1.049628776243469
$ws.Range("L6").Value = 1.048876995907549
$ws.Range("M6").Value = 1.058939924073876
$ws.Range("N6").Value = 1.019574738040396

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042512642519671
$ws.Range("D7").Value = 1.046927693548218
$ws.Range("E7").Value = 1.046085545678461
$ws.Range("F7").Value = 1.056215126907682
$ws.Range("I7").Value = 1.037750276454477
$ws.Range("J7").Value = 1.047011246139385
$ws.Range("K7").Value = 1.049383650646342
$ws.Range("L7").Value = 1.048543586090756
$ws.Range("M7").Value = 1.058648345576949
$ws.Range("N7").Value = 1.019473250461193

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04080168484398
$ws.Range("D8").Value = 1.045651886796342
$ws.Range("E8").Value = 1.044440592398926
$ws.Range("F8").Value = 1.054747715074198
$ws.Range("I8").Value = 1.037408334720537
$ws.Range("J8").Value = 1.045768430166127
$ws.Range("K8").Value = 1.048356702870511
$ws.Range("L8").Value = 1.047148730129633
$ws.Range("M8").Value = 1.057427851979467
$ws.Range("N8").Value = 1.01904792480516

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037783287892882
$ws.Range("D9").Value = 1.043399459913415
$ws.Range("E9").Value = 1.041541324112518
$ws.Range("F9").Value = 1.052159761454409
$ws.Range("I9").Value = 1.036793683662754
$ws.Range("J9").Value = 1.043571451920313
$ws.Range("K9").Value = 1.04653808714449
$ws.Range("L9").Value = 1.044685930710315
$ws.Range("M9").Value = 1.055270509724927
$ws.Range("N9").Value = 1.018294172452663

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035768485714018
$ws.Range("D10").Value = 1.041894863266774
$ws.Range("E10").Value = 1.03960788397022
$ws.Range("F10").Value = 1.050432878596051
$ws.Range("I10").Value = 1.036375743017186
$ws.Range("J10").Value = 1.042101982143152
$ws.Range("K10").Value = 1.045319531562788
$ws.Range("L10").Value = 1.043040641333361
$ws.Range("M10").Value = 1.053827690541845
$ws.Range("N10").Value = 1.017788761611106

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034895406270788
$ws.Range("D11").Value = 1.041242628320379
$ws.Range("E11").Value = 1.03877050627312
$ws.Range("F11").Value = 1.049684718977425
$ws.Range("I11").Value = 1.03619282357129
$ws.Range("J11").Value = 1.04146451408124
$ws.Range("K11").Value = 1.044790404771682
$ws.Range("L11").Value = 1.042327369396024
$ws.Range("M11").Value = 1.053201821500687
$ws.Range("N11").Value = 1.01756921486299

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.03457100296521
$ws.Range("D12").Value = 1.041000246893529
$ws.Range("E12").Value = 1.038459435893605
$ws.Range("F12").Value = 1.049406755601509
$ws.Range("I12").Value = 1.036124585673831
$ws.Range("J12").Value = 1.041227550517587
$ws.Range("K12").Value = 1.044593638617037
$ws.Range("L12").Value = 1.042062297808667
$ws.Range("M12").Value = 1.05296917566062
$ws.Range("N12").Value = 1.0174875594283

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.03464059328924
$ws.Range("D13").Value = 1.041052243646868
$ws.Range("E13").Value = 1.038526163000685
$ws.Range("F13").Value = 1.049466382603372
$ws.Range("I13").Value = 1.036139236226692
$ws.Range("J13").Value = 1.041278388172411
$ws.Range("K13").Value = 1.04463585583357
$ws.Range("L13").Value = 1.042119162520585
$ws.Range("M13").Value = 1.053019086719674
$ws.Range("N13").Value = 1.01750507960495

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034868593111414
$ws.Range("D14").Value = 1.041222595298791
$ws.Range("E14").Value = 1.038744793749669
$ws.Range("F14").Value = 1.049661743748415
$ws.Range("I14").Value = 1.036187188996452
$ws.Range("J14").Value = 1.041444930276651
$ws.Range("K14").Value = 1.044774144619115
$ws.Range("L14").Value = 1.04230546117707
$ws.Range("M14").Value = 1.05318259440755
$ws.Range("N14").Value = 1.017562467365963

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03500905759316
$ws.Range("D15").Value = 1.041327539630857
$ws.Range("E15").Value = 1.038879495164492
$ws.Range("F15").Value = 1.049782103703609
$ws.Range("I15").Value = 1.036216695370388
$ws.Range("J15").Value = 1.041547518495541
$ws.Range("K15").Value = 1.044859319044176
$ws.Range("L15").Value = 1.042420228533719
$ws.Range("M15").Value = 1.053283314272143
$ws.Range("N15").Value = 1.017597811791327

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035826414937871
$ws.Range("D16").Value = 1.041938134307178
$ws.Range("E16").Value = 1.03966345379456
$ws.Range("F16").Value = 1.050482522695609
$ws.Range("I16").Value = 1.036387841664373
$ws.Range("J16").Value = 1.042144263761941
$ws.Range("K16").Value = 1.045354616531322
$ws.Range("L16").Value = 1.043087960678111
$ws.Range("M16").Value = 1.053869203655861
$ws.Range("N16").Value = 1.017803317382514

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036338942906365
$ws.Range("D17").Value = 1.042320946149929
$ws.Range("E17").Value = 1.040155158471306
$ws.Range("F17").Value = 1.050921766190934
$ws.Range("I17").Value = 1.036494675047506
$ws.Range("J17").Value = 1.042518269120808
$ws.Range("K17").Value = 1.045664904846563
$ws.Range("L17").Value = 1.043506581651837
$ws.Range("M17").Value = 1.054236415554004
$ws.Range("N17").Value = 1.017932037532169

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036637828702324
$ws.Range("D18").Value = 1.042544163069709
$ws.Range("E18").Value = 1.040441943822037
$ws.Range("F18").Value = 1.051177930286237
$ws.Range("I18").Value = 1.036556801212355
$ws.Range("J18").Value = 1.042736306638774
$ws.Range("K18").Value = 1.045845747700411
$ws.Range("L18").Value = 1.043750674286943
$ws.Range("M18").Value = 1.054450496237138
$ws.Range("N18").Value = 1.018007050279272

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036739730454376
$ws.Range("D19").Value = 1.042620262328654
$ws.Range("E19").Value = 1.040539727313863
$ws.Range("F19").Value = 1.051265269041863
$ws.Range("I19").Value = 1.036577952773555
$ws.Range("J19").Value = 1.042810632626733
$ws.Range("K19").Value = 1.045907386230557
$ws.Range("L19").Value = 1.043833889761001
$ws.Range("M19").Value = 1.054523473995881
$ws.Range("N19").Value = 1.018032616254366

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036283960062941
$ws.Range("D20").Value = 1.042279881395343
$ws.Range("E20").Value = 1.040102405064636
$ws.Range("F20").Value = 1.050874643593591
$ws.Range("I20").Value = 1.036483232281892
$ws.Range("J20").Value = 1.04247815364276
$ws.Range("K20").Value = 1.045631628660154
$ws.Range("L20").Value = 1.043461676075533
$ws.Range("M20").Value = 1.054197028359638
$ws.Range("N20").Value = 1.017918234069925

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034801455714878
$ws.Range("D21").Value = 1.04117243409971
$ws.Range("E21").Value = 1.03868041332512
$ws.Range("F21").Value = 1.049604216534093
$ws.Range("I21").Value = 1.036173076208202
$ws.Range("J21").Value = 1.041395892754834
$ws.Range("K21").Value = 1.044733428242509
$ws.Range("L21").Value = 1.042250604494815
$ws.Range("M21").Value = 1.053134450184561
$ws.Range("N21").Value = 1.017545571032801

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033868750031116
$ws.Range("D22").Value = 1.040475487165951
$ws.Range("E22").Value = 1.037786169635742
$ws.Range("F22").Value = 1.048805079609683
$ws.Range("I22").Value = 1.03597637047014
$ws.Range("J22").Value = 1.040714392323683
$ws.Range("K22").Value = 1.044167391799114
$ws.Range("L22").Value = 1.041488398530488
$ws.Range("M22").Value = 1.052465379262496
$ws.Range("N22").Value = 1.017310649347741

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034363252844446
$ws.Range("D23").Value = 1.040845014270695
$ws.Range("E23").Value = 1.038260243357525
$ws.Range("F23").Value = 1.049228752959367
$ws.Range("I23").Value = 1.036080809133695
$ws.Range("J23").Value = 1.041075768041897
$ws.Range("K23").Value = 1.044467582504567
$ws.Range("L23").Value = 1.041892531048095
$ws.Range("M23").Value = 1.052820160575669
$ws.Range("N23").Value = 1.017435244196262

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03630880463798
$ws.Range("D24").Value = 1.042298437008958
$ws.Range("E24").Value = 1.040126242112893
$ws.Range("F24").Value = 1.050895936388107
$ws.Range("I24").Value = 1.036488403355836
$ws.Range("J24").Value = 1.042496280451118
$ws.Range("K24").Value = 1.045646665178485
$ws.Range("L24").Value = 1.043481967225511
$ws.Range("M24").Value = 1.05421482607046
$ws.Range("N24").Value = 1.017924471469035

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038564048128122
$ws.Range("D25").Value = 1.043982284433217
$ws.Range("E25").Value = 1.042290947694295
$ws.Range("F25").Value = 1.052829079930753
$ws.Range("I25").Value = 1.036954023793399
$ws.Range("J25").Value = 1.04414026330906
$ws.Range("K25").Value = 1.047009319217289
$ws.Range("L25").Value = 1.0453232165328
$ws.Range("M25").Value = 1.055829036089069
$ws.Range("N25").Value = 1.018489545660981
